$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New cell values (order chosen to match shared-string insertion order)
$ws.Range("A3").Value = "rutaCRM"
$ws.Range("B3").Value = "C:\Users\Gonzalo PL\Desktop\Automatizaciones\Automatizacione-EMP1"

$ws.Range("A4").Value = "LISTA PRECIOS"
$ws.Range("A5").Value = "LISTA PRECIOS LG"

$ws.Range("B5").Value = "C:\Users\Gonzalo PL\Desktop\Automatizaciones\Automatizacione-EMP1\Codigos Costos\LISTA DE PRECIOS LG - VENTAS COSTOS 30.10.25- OCTUBRE 2025.xlsx"
$ws.Range("B4").Value = "C:\Users\Gonzalo PL\Desktop\Automatizaciones\Automatizacione-EMP1\Codigos Costos\LISTA DE PRECIOS - VENTAS COSTOS 20.10.25- OCTUBRE 2025-4.xlsx"

# Styles: borders first on the whole A3:B5 block
$full = $ws.Range("A3:B5")
$full.Borders.LineStyle = 1

# A4:A5 fill (theme accent2, lighter 60%)
$ws.Range("A4:A5").Interior.ThemeColor = 6

# A3 fill (yellow)
$ws.Range("A3").Interior.Color = 65535
